$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to set a cell as plain text, avoiding Excel auto-converting
# numeric-looking strings (e.g. "1.00", "116.99") into numbers, while
# keeping the cell style identical to the original (no explicit style).
function Set-TextValue {
    param($Ref, $Value)
    $cell = $ws.Range($Ref)
    $cell.NumberFormat = "@"
    $cell.Value = $Value
    $cell.Style = "Normal"
}

Set-TextValue "D2" "62.373.61"
Set-TextValue "E2" "  +9.33%  "
Set-TextValue "D3" "3.375.21"
Set-TextValue "E3" "  +3.99%  "
Set-TextValue "E4" "  +0.01%  "
Set-TextValue "E5" "  +4.81%  "
Set-TextValue "D6" "116.99"
Set-TextValue "E6" "  +8.15%  "
Set-TextValue "D7" "3.367.02"
Set-TextValue "E7" "  +3.85%  "
Set-TextValue "E8" "  -2.10%  "
Set-TextValue "D9" "1.00"
Set-TextValue "E9" "  -0.01%  "
Set-TextValue "D10" "0.631"
Set-TextValue "E10" "  +1.07%  "
Set-TextValue "D11" "0.116"
Set-TextValue "E11" "  +18.13%  "
Set-TextValue "D12" "40.03"
Set-TextValue "E12" "  +2.02%  "
Set-TextValue "E13" "  -0.51%  "
Set-TextValue "D14" "3.900.47"
Set-TextValue "E14" "  +3.74%  "
Set-TextValue "D15" "8.37"
Set-TextValue "E15" "  -0.10%  "
Set-TextValue "E16" "  +1.57%  "
Set-TextValue "D17" "3.390.99"
Set-TextValue "E17" "  +0.95%  "
Set-TextValue "D18" "62.201.43"
Set-TextValue "E18" "  +9.21%  "
Set-TextValue "E19" "  -1.87%  "
Set-TextValue "D20" "10.95"
Set-TextValue "E20" "  +1.25%  "
Set-TextValue "D21" "0.0000117"
Set-TextValue "E21" "  +7.54%  "
Set-TextValue "E22" "  +0.60%  "
Set-TextValue "D23" "12.62"
Set-TextValue "E23" "  -3.51%  "
Set-TextValue "D24" "297.54"
Set-TextValue "E24" "  +1.25%  "
Set-TextValue "D25" "74.89"
Set-TextValue "E25" "  +1.01%  "
Set-TextValue "E26" "  -0.69%  "
Set-TextValue "D27" "29.67"
Set-TextValue "E27" "  +5.63%  "
Set-TextValue "D28" "7.96"
Set-TextValue "E28" "  +10.52%  "
Set-TextValue "D29" "0.177"
Set-TextValue "E29" "  +4.69%  "
Set-TextValue "D30" "4.26"
Set-TextValue "E30" "  -2.16%  "
Set-TextValue "D31" "7.61"
Set-TextValue "E31" "  -0.45%  "
Set-TextValue "D32" "43.21"
Set-TextValue "E32" "  +7.28%  "
Set-TextValue "E33" "  +4.35%  "
Set-TextValue "D34" "11.45"
Set-TextValue "E34" "  +2.12%  "
Set-TextValue "E35" "  +19.89%  "
Set-TextValue "D36" "1.00"
Set-TextValue "E36" "  +0.02%  "
Set-TextValue "D37" "0.0491"
Set-TextValue "E37" "  +0.43%  "
Set-TextValue "D38" "52.28"
Set-TextValue "E38" "  +1.45%  "
Set-TextValue "D39" "3.12"
Set-TextValue "E39" "  +5.86%  "
Set-TextValue "D40" "0.999"
Set-TextValue "E40" "  -0.04%  "
Set-TextValue "D41" "3.44"
Set-TextValue "E41" "  -0.72%  "
Set-TextValue "D42" "133.40"
Set-TextValue "E42" "  -4.28%  "
Set-TextValue "E43" "  -1.24%  "
Set-TextValue "B44" "TheGraph"
Set-TextValue "C44" "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D44" "0.289"
Set-TextValue "E44" "  +3.62%  "
Set-TextValue "B45" "ARBITRUM"
Set-TextValue "C45" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D45" "1.90"
Set-TextValue "E45" "  +0.50%  "
Set-TextValue "E46" "  -1.91%  "
Set-TextValue "D47" "16.50"
Set-TextValue "E47" "  -3.17%  "
Set-TextValue "E48" "  -4.18%  "
Set-TextValue "B49" "EnergySwap"
Set-TextValue "C49" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D49" "21.26"
Set-TextValue "E49" "  -4.10%  "
Set-TextValue "B50" "Maker"
Set-TextValue "C50" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D50" "2.172.24"
Set-TextValue "E50" "  +0.62%  "
Set-TextValue "D51" "3.701.85"
Set-TextValue "E51" "  +3.67%  "
